$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell values in column A (rows 1-33) ---
$ws.Cells.Item(1, 1).Value = 0.09472458870376954
$ws.Cells.Item(2, 1).Value = -0.005999999967798431
$ws.Cells.Item(3, 1).Value = -0.003999999975107471
$ws.Cells.Item(4, 1).Value = -0.007999999952783554
$ws.Cells.Item(5, 1).Value = -0.0029999999788667964
$ws.Cells.Item(6, 1).Value = -0.001999999982777112
$ws.Cells.Item(7, 1).Value = -0.009999999939823923
$ws.Cells.Item(8, 1).Value = -0.009999999940711657
$ws.Cells.Item(9, 1).Value = -0.001999999985863532
$ws.Cells.Item(10, 1).Value = 0.054054599954188376
$ws.Cells.Item(11, 1).Value = -0.002999999982549184
$ws.Cells.Item(12, 1).Value = -0.0034999999793075
$ws.Cells.Item(13, 1).Value = -0.003499999975863588
$ws.Cells.Item(14, 1).Value = -0.007999999950761172
$ws.Cells.Item(15, 1).Value = -0.000999999987759459
$ws.Cells.Item(16, 1).Value = -0.001999999981836975
$ws.Cells.Item(17, 1).Value = -0.0019999999810016433
$ws.Cells.Item(18, 1).Value = -0.003999999970089263
$ws.Cells.Item(19, 1).Value = -0.003999999978431035
$ws.Cells.Item(20, 1).Value = -0.003999999976571189
$ws.Cells.Item(21, 1).Value = -0.003999999976282531
$ws.Cells.Item(22, 1).Value = -0.003999999976108448
$ws.Cells.Item(23, 1).Value = -0.004999999967937541
$ws.Cells.Item(24, 1).Value = -0.01999999988375034
$ws.Cells.Item(25, 1).Value = -0.01999999988204859
$ws.Cells.Item(26, 1).Value = -0.0024999999799071304
$ws.Cells.Item(27, 1).Value = -0.0024999999796349037
$ws.Cells.Item(28, 1).Value = -0.001999999981308065
$ws.Cells.Item(29, 1).Value = 0.053019046948836746
$ws.Cells.Item(30, 1).Value = -0.035378971998003106
$ws.Cells.Item(31, 1).Value = -0.006999999948879676
$ws.Cells.Item(32, 1).Value = -0.009999999932393422
$ws.Cells.Item(33, 1).Value = -0.003999999964495515

# --- Widen column A from 15.42578125 to 16.42578125 characters ---
# The COM ColumnWidth setter here quantizes to an internal pixel grid
# (steps of 1/6 character-width units), so we pick the input that lands
# on the grid point closest to the target stored width (16.5).
$ws.Range("A1").EntireColumn.ColumnWidth = 15.666666666666668
